$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("file title JSON")

$rows = @(
    @(678, '01ride2015.mp3', 'ride01'),
    @(679, '02ride2015.mp3', 'ride02'),
    @(680, '03%20YYZ.mp3', 'rush xyz'),
    @(681, '03ride2015.mp3', 'ride03'),
    @(682, '04ride2015.mp3', 'ride04'),
    @(683, '05ride2015.mp3', 'ride05'),
    @(684, '06ride2015.mp3', 'ride06'),
    @(685, '07ride2015.mp3', 'ride07'),
    @(686, '08ride2015.mp3', 'ride08'),
    @(687, '09ride2015.mp3', 'ride09'),
    @(688, '10ride2015.mp3', 'ride10'),
    @(689, '11ride2015.mp3', 'ride11'),
    @(690, '12ride2015.mp3', 'ride12'),
    @(691, '13ride2015.mp3', 'ride13'),
    @(692, '14ride2015.mp3', 'ride14'),
    @(693, 'ZOOM0012-SkerikTrio-OwlnThistle-20130227.mp3', 'Skerik Trio-Owl n Thistle'),
    @(694, 'blondie-callme-7inch.mp3', 'Blondie-Call Me'),
    @(695, '01 Here To Fall.mp3', 'Yo La Tengo-Here To Fall'),
    @(696, '01ride2015.mp3', 'Ride01'),
    @(697, '02 Avalon Or Someone Very Similar.mp3', 'Yo La Tengo-Avalon Or Someone Very Similar'),
    @(698, '02ride2015.mp3', 'Ride02'),
    @(699, '03 By Twos.mp3', 'Yo La Tengo-By Twos'),
    @(700, '03 YYZ.mp3', 'Rush YYZ'),
    @(701, '03ride2015.mp3', 'Ride03'),
    @(702, '04 Nothing To Hide.mp3', 'Yo La Tengo-Nothing To Hide'),
    @(703, '04ride2015.mp3', 'Ride04'),
    @(704, '05 Periodically Triple Or Double.mp3', 'Yo La Tengo-Periodically Triple Or Double'),
    @(705, '05ride2015.mp3', 'Ride05'),
    @(706, '06 If Its True.mp3', 'Yo La Tengo-If Its True'),
    @(707, '06ride2015.mp3', 'Ride 06'),
    @(708, '07 Im On My Way.mp3', 'Yo La Tengo-Im On My Way'),
    @(709, '07ride2015.mp3', 'Ride07'),
    @(710, '08 When Its Dark.mp3', 'Yo La Tengo-When Its Dark'),
    @(711, '08ride2015.mp3', 'Ride08'),
    @(712, '09 All Your Secrets.mp3', 'Yo La Tengo-All Your Secrets'),
    @(713, '09ride2015.mp3', 'Ride09'),
    @(714, '10 More Stars Than There Are In Heaven.mp3', 'Yo La Tengo-More Stars Than There Are In Heaven'),
    @(715, '10ride2015.mp3', 'Ride10'),
    @(716, '11 The Fireside.mp3', 'Yo La Tengo-The Fireside'),
    @(717, '11ride2015.mp3', 'Ride11'),
    @(718, '12ride2015.mp3', 'Ride12'),
    @(719, '13ride2015.mp3', 'Ride13'),
    @(720, '14ride2015.mp3', 'Ride14'),
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $fileName = $r[1]
    $title = $r[2]
    $ws.Cells.Item($rowNum, 1).Value2 = $fileName
    $ws.Cells.Item($rowNum, 2).Value2 = $title
    $aRef = "A" + $rowNum
    $bRef = "B" + $rowNum
    $ws.Cells.Item($rowNum, 3).Formula = "=""{'file':'""&$aRef&""','title':'""&$bRef&""'},"""
}

$ws.Activate()
$ws.Range("C693:C720").Select()

Write-Host "done applying rows"
